$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -1
$ws.Range("B1").Value = 3.43013596534729
$ws.Range("C1").Value = 1.90750777721405
$ws.Range("D1").Value = 1.44348156452179
$ws.Range("E1").Value = 1.282873749732971
